$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 32 (pushes strDlgReset and everything below it down by one),
# copying formatting from the surrounding rows automatically.
$ws.Rows("32:32").Insert()

# Populate the newly inserted row with the new "strWindowPos" translation entry.
$ws.Range("B32").Value = "localization\strings"
$ws.Range("C32").Value = "strWindowPos"
$ws.Range("D32").Value = 'In "settings" form, tab "User interface"'
$ws.Range("E32").Value = "Remember window position and size on startup"

# Row 25 (strChkDlgPath) gets the same comment text and grows to a 2-line row.
$ws.Range("D25").Value = 'In "settings" form, tab "User interface"'
$ws.Rows("25:25").RowHeight = 30

# Keep the "Tabla13" table in sync with the newly inserted row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B2:F204"))
